$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column U (shifts old U -> V)
$ws.Columns("U:U").Insert()

# New column U header (no special style) + MD5/version data
$ws.Range("U1").Value = "MD5 & Script Version"
$ws.Range("U1").Style = "Normal"
$ws.Range("U2").Value = "MD5: 33d07e47ee89e1db338f01228abafde5 | Script: v3.0.0"
$ws.Range("U3").Value = "MD5: 2363d6688c56d6cd78be8f9a481fcf8f | Script: v3.0.0"
$ws.Range("U4").Value = "MD5: 21b9a9f0d70e550ae3d66ed506ff9e2d | Script: v3.0.0"
$ws.Range("U5").Value = "MD5: 305f3681692ea4f39d05961cd6714eeb | Script: v3.0.0"
